# Soulreaping weapons cannot be disenchanted
#
# Adds a new "Soulreaping" weapon-material row to the Weapons sheet
# (sorted alphabetically between "SkyforgeSteel" and "Spectral"), which
# pushes it into row 34 and shifts every row below it down by one.
# The new entry copies the shape of similar low-tier "Craftsmanship"
# materials (e.g. AncientNordHoned) - Damage 3, Weight 1, Gold 2/3,
# tempered/broken down with Steel Ingot, perk Craftsmanship - but is
# intentionally left without a Breakdown/disenchant recipe so it cannot
# be disenchanted.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weapons")
$ws4 = $wb.Worksheets.Item("Artifacts")

# --- Insert the new row ------------------------------------------------
$ws1.Rows.Item(34).Insert()

$ws1.Cells.Item(34, 1).Value2 = "Soulreaping"
$ws1.Cells.Item(34, 2).Value2 = 3
$ws1.Cells.Item(34, 3).Value2 = 1
$ws1.Cells.Item(34, 4).Value2 = 0.66666666666666663
$ws1.Cells.Item(34, 4).NumberFormat = "# ?/?"
$ws1.Cells.Item(34, 7).Value2 = "Steel Ingot"
$ws1.Cells.Item(34, 8).Value2 = "Steel Ingot"
$ws1.Cells.Item(34, 9).Value2 = "Craftsmanship"

# --- Restore/update the window & selection state ------------------------
# Weapons becomes the active tab, scrolled near the newly added row, with
# A39 selected.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A39").Select()

# Artifacts keeps its own scroll position/selection but is no longer the
# active tab.
$ws4.Select()
$excel.ActiveWindow.ScrollRow = 51
$excel.ActiveWindow.ScrollColumn = 1
$ws4.Range("C71").Select()

$ws1.Activate()
